$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("104_1")

$ws.Range("B11").Value = 1
$ws.Range("B33").Value = 1

$wb.Save()
